$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Move the existing "note" row (old row 102: A102 empty / B102 footnote text)
#        down to row 103, preserving its values & formatting exactly. ---
$ws.Range("A102:B102").Copy($ws.Range("A103:B103"))

# --- 2. Turn (old, now vacated) row 102 into a fresh data row matching the
#        style of the row above it (row 101), then fill in the new data. ---
$ws.Range("A101:E101").Copy($ws.Range("A102:E102"))

$ws.Range("A102").Value = 43957
$ws.Range("B102").Value = 455
$ws.Range("C102").Value = 34240
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 6958

$excel.CutCopyMode = $false

# --- 3. Update the workbook-level print area so it now spans through the
#        new last row (104 = header + 102 data rows + 1 footnote row). ---
$wb.Names.Item(1).RefersTo = '=相談件数!$A$1:$E$104'

# --- 4. Update the view state: scroll the frozen pane down a bit and move
#        the active selection onto the new footnote row. ---
$excel.ActiveWindow.ScrollRow = 89
$ws.Range("C103").Select() | Out-Null
